$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 3
$ws.Range("H3").Value = 50000
$ws.Range("J3").Value = 50000
$ws.Range("L3").Value = 50000
$ws.Range("N3").Value = -50228
# Row 17
$ws.Range("H17").Value = 4538.0605
$ws.Range("J17").Value = 4791.8667
$ws.Range("L17").Value = 14375.6001
$ws.Range("N17").Value = -14711.6001
# Row 43
$ws.Range("H43").Value = 3124.125
$ws.Range("J43").Value = 3959.6
$ws.Range("L43").Value = 3959.6
$ws.Range("N43").Value = -4097.6
# Row 100
$ws.Range("H100").Value = 2522.7144
$ws.Range("I100").Value = 1168.3334
$ws.Range("J100").Value = 3538.5
$ws.Range("K100").Value = 1168.3334
$ws.Range("L100").Value = 3538.5
$ws.Range("M100").Value = -627.3334
$ws.Range("N100").Value = -4620.5
# Row 102
$ws.Range("H102").Value = 50000
$ws.Range("J102").Value = 50000
$ws.Range("L102").Value = 50000
$ws.Range("N102").Value = -56490
# Row 105
$ws.Range("H105").Value = 35000
$ws.Range("J105").Value = 35000
$ws.Range("L105").Value = 35000
$ws.Range("N105").Value = -41988
# Row 110
$ws.Range("H110").Value = 50000
$ws.Range("J110").Value = 50000
$ws.Range("L110").Value = 50000
$ws.Range("N110").Value = -58180
# Row 129
$ws.Range("H129").Value = 678.9091
$ws.Range("I129").Value = 678.9091
$ws.Range("K129").Value = 2036.7273
$ws.Range("M129").Value = 2963.2727
# Row 137
$ws.Range("H137").Value = 1788.2142
$ws.Range("I137").Value = 1653.7
$ws.Range("K137").Value = 4961.1
$ws.Range("M137").Value = -2411.1

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 5
$ws.Range("H5").Value = 171.8
$ws.Range("I5").Value = 183.33333
$ws.Range("K5").Value = 183.33333
$ws.Range("M5").Value = -71.33332999999999
# Row 61
$ws.Range("H61").Value = 4139.7144
$ws.Range("I61").Value = 3092.25
$ws.Range("K61").Value = 3092.25
$ws.Range("M61").Value = -2880.25
# Row 74
$ws.Range("H74").Value = 5219.5293
$ws.Range("I74").Value = 3417.5454
$ws.Range("K74").Value = 3417.5454
$ws.Range("M74").Value = -2543.5454
# Row 77
$ws.Range("H77").Value = 5219.5293
$ws.Range("I77").Value = 3417.5454
$ws.Range("K77").Value = 17087.727
$ws.Range("M77").Value = -12719.727
# Row 132
$ws.Range("H132").Value = 1424.5306
$ws.Range("I132").Value = 1014.175
$ws.Range("K132").Value = 3042.525
$ws.Range("M132").Value = -512.5249999999996
# Row 136
$ws.Range("H136").Value = 4139.7144
$ws.Range("I136").Value = 3092.25
$ws.Range("K136").Value = 9276.75
$ws.Range("M136").Value = -6726.75

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 4
$ws.Range("H4").Value = 171.8
$ws.Range("I4").Value = 183.33333
$ws.Range("K4").Value = 183.33333
$ws.Range("M4").Value = -68.33332999999999
# Row 44
$ws.Range("H44").Value = 23163.334
$ws.Range("J44").Value = 23163.334
$ws.Range("L44").Value = 23163.334
$ws.Range("N44").Value = -24157.334

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 9
$ws.Range("H9").Value = 67885
$ws.Range("J9").Value = 67885
$ws.Range("L9").Value = 67885
$ws.Range("N9").Value = -68221
# Row 31
$ws.Range("H31").Value = 7453.951
$ws.Range("I31").Value = 3905.45
$ws.Range("K31").Value = 3905.45
$ws.Range("M31").Value = -3610.45
# Row 34
$ws.Range("H34").Value = 7453.951
$ws.Range("I34").Value = 3905.45
$ws.Range("K34").Value = 3905.45
$ws.Range("M34").Value = -3703.45
# Row 43
$ws.Range("H43").Value = 10332.667
$ws.Range("J43").Value = 10332.667
$ws.Range("L43").Value = 10332.667
$ws.Range("N43").Value = -10700.667
# Row 58
$ws.Range("H58").Value = 2290.516
$ws.Range("I58").Value = 1865.2778
$ws.Range("J58").Value = 2879.3076
$ws.Range("K58").Value = 1865.2778
$ws.Range("L58").Value = 2879.3076
$ws.Range("M58").Value = -1662.2778
$ws.Range("N58").Value = -3285.3076
# Row 101
$ws.Range("H101").Value = 10332.667
$ws.Range("J101").Value = 10332.667
$ws.Range("L101").Value = 10332.667
$ws.Range("N101").Value = -16822.667
# Row 132
$ws.Range("H132").Value = 2272.7666
$ws.Range("I132").Value = 1384.463
$ws.Range("J132").Value = 10267.5
$ws.Range("K132").Value = 4153.389
$ws.Range("L132").Value = 30802.5
$ws.Range("M132").Value = -1623.389
$ws.Range("N132").Value = -35862.5
# Row 136
$ws.Range("H136").Value = 2290.516
$ws.Range("I136").Value = 1865.2778
$ws.Range("J136").Value = 2879.3076
$ws.Range("K136").Value = 5595.8334
$ws.Range("L136").Value = 8637.9228
$ws.Range("M136").Value = -3045.8334
$ws.Range("N136").Value = -13737.9228

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 1356
$ws.Range("I5").Value = 1193.5
$ws.Range("J5").Value = 1464.3334
$ws.Range("K5").Value = 3580.5
$ws.Range("L5").Value = 4393.0002
$ws.Range("M5").Value = -3468.5
$ws.Range("N5").Value = -4617.0002
# Row 104
$ws.Range("H104").Value = 0
$ws.Range("J104").Value = 0
$ws.Range("L104").Value = 0
$ws.Range("N104").ClearContents()
# Row 107
$ws.Range("H107").Value = 441.36365
$ws.Range("I107").Value = 348
$ws.Range("J107").Value = 476.375
$ws.Range("K107").Value = 1044
$ws.Range("L107").Value = 1429.125
$ws.Range("M107").Value = 876
$ws.Range("N107").Value = -5269.125
# Row 122
$ws.Range("H122").Value = 939.8
$ws.Range("I122").Value = 1500
$ws.Range("J122").Value = 799.75
$ws.Range("K122").Value = 13500
$ws.Range("L122").Value = 7197.75
$ws.Range("M122").Value = -11050
$ws.Range("N122").Value = -12097.75
# Row 135
$ws.Range("H135").Value = 1356
$ws.Range("I135").Value = 1193.5
$ws.Range("J135").Value = 1464.3334
$ws.Range("K135").Value = 10741.5
$ws.Range("L135").Value = 13179.0006
$ws.Range("M135").Value = -8206.5
$ws.Range("N135").Value = -18249.0006

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 100
$ws.Range("H100").Value = 111555
$ws.Range("J100").Value = 111555
$ws.Range("L100").Value = 111555
$ws.Range("N100").Value = -113719

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 7487.7646
$ws.Range("I7").Value = 6163.7144
$ws.Range("J7").Value = 13666.667
$ws.Range("K7").Value = 6163.7144
$ws.Range("L7").Value = 13666.667
$ws.Range("M7").Value = -6051.7144
$ws.Range("N7").Value = -13890.667
# Row 46
$ws.Range("H46").Value = 1908.3334
$ws.Range("J46").Value = 2342.5
$ws.Range("L46").Value = 2342.5
$ws.Range("N46").Value = -2718.5
# Row 98
$ws.Range("H98").Value = 0
$ws.Range("J98").Value = 0
$ws.Range("L98").Value = 0
$ws.Range("N98").ClearContents()
# Row 126
$ws.Range("H126").Value = 7487.7646
$ws.Range("I126").Value = 6163.7144
$ws.Range("J126").Value = 13666.667
$ws.Range("K126").Value = 18491.1432
$ws.Range("L126").Value = 41000.001
$ws.Range("M126").Value = -16021.1432
$ws.Range("N126").Value = -45940.001
# Row 136
$ws.Range("H136").Value = 6590.844
$ws.Range("J136").Value = 9469.549999999999
$ws.Range("L136").Value = 28408.65
$ws.Range("N136").Value = -33508.64999999999

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 96
$ws.Range("H96").Value = 10028.143
$ws.Range("J96").Value = 12059.4
$ws.Range("L96").Value = 12059.4
$ws.Range("N96").Value = -14805.4
# Row 104
$ws.Range("H104").Value = 29249.75
$ws.Range("J104").Value = 29249.75
$ws.Range("L104").Value = 29249.75
$ws.Range("N104").Value = -36237.75
# Row 107
$ws.Range("H107").Value = 1422.8108
$ws.Range("I107").Value = 1631.5454
$ws.Range("J107").Value = 1116.6666
$ws.Range("K107").Value = 4894.6362
$ws.Range("L107").Value = 3349.9998
$ws.Range("M107").Value = -2974.6362
$ws.Range("N107").Value = -7189.9998
# Row 110
$ws.Range("H110").Value = 150000
$ws.Range("J110").Value = 150000
$ws.Range("L110").Value = 150000
$ws.Range("N110").Value = -158180
# Row 132
$ws.Range("H132").Value = 1755.356
$ws.Range("I132").Value = 1382.7234
$ws.Range("K132").Value = 4148.1702
$ws.Range("M132").Value = -1618.1702
